$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (E1:G1)
$ws.Range("E1").Value = "Qtd Vendas"
$ws.Range("F1").Value = "Ticket Médio Anual"
$ws.Range("G1").Value = "Evolução Ticket Médio (%)"

# Apply the same style as the existing header cells (A1:D1) to the new headers
$ws.Range("A1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

# New data values for E (Qtd Vendas), F (Ticket Médio Anual), G (Evolução Ticket Médio %)
$ws.Range("E2").Value = 258
$ws.Range("F2").Value = 448.3624418604651

$ws.Range("E3").Value = 2141
$ws.Range("F3").Value = 445.1697244278375
$ws.Range("G3").Value = -0.712084049542494

$ws.Range("E4").Value = 2716
$ws.Range("F4").Value = 657.8443262150221
$ws.Range("G4").Value = 47.77382425557543

$ws.Range("E5").Value = 3342
$ws.Range("F5").Value = 863.5469838420108
$ws.Range("G5").Value = 31.26919993526753

$ws.Range("E6").Value = 4662
$ws.Range("F6").Value = 968.9902981552981
$ws.Range("G6").Value = 12.21048956064428

$ws.Range("E7").Value = 1887
$ws.Range("F7").Value = 1046.919337572867
$ws.Range("G7").Value = 8.04229305142945
